{"js": "// Georgia -> DejaVu Sans for the styles that reference Georgia in this\n// document's style sheet (sw/qa/core/layout/data/table-fly-overlap-spacing.docx\n// regression fixture): Normal, Heading 1 Char, Heading 2 Char, Header Char,\n// Footer Char, Title Char, No Spacing, Subtitle Char, Body Text Char,\n// No Spacing Char, Comment Char.\n\nconst OLD_FONT = \"Georgia\";\nconst NEW_FONT = \"DejaVu Sans\";\n\nconst styleNames = [\n  \"Normal\",\n  \"Heading 1 Char\",\n  \"Heading 2 Char\",\n  \"Header Char\",\n  \"Footer Char\",\n  \"Title Char\",\n  \"No Spacing\",\n  \"Subtitle Char\",\n  \"Body Text Char\",\n  \"No Spacing Char\",\n  \"Comment Char\",\n];\n\nconst styles = context.document.getStyles();\n\n// Resolve every style up front, then load its font name so we only flip the\n// ones that are actually still set to the old font (defensive / idempotent).\nconst styleObjs = styleNames.map((name) => styles.getByNameOrNullObject(name));\nstyleObjs.forEach((s) => s.load(\"isNullObject,nameLocal\"));\nawait context.sync();\n\nconst fonts = [];\nfor (const s of styleObjs) {\n  if (!s.isNullObject) {\n    s.font.load(\"name\");\n    fonts.push(s.font);\n  } else {\n    fonts.push(null);\n  }\n}\nawait context.sync();\n\nfor (let i = 0; i < styleObjs.length; i++) {\n  const s = styleObjs[i];\n  const font = fonts[i];\n  if (s.isNullObject || !font) {\n    continue;\n  }\n  if (font.name === OLD_FONT) {\n    font.name = NEW_FONT;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Georgia -> DejaVu Sans for the styles that reference Georgia in this\n# document's style sheet (sw/qa/core/layout/data/table-fly-overlap-spacing.docx\n# regression fixture): Normal, Heading 1 Char, Heading 2 Char, Header Char,\n# Footer Char, Title Char, No Spacing, Subtitle Char, Body Text Char,\n# No Spacing Char, Comment Char.\n\n$d = $word.ActiveDocument\n\n$OldFont = \"Georgia\"\n$NewFont = \"DejaVu Sans\"\n\n$styleNames = @(\n    \"Normal\",\n    \"Heading 1 Char\",\n    \"Heading 2 Char\",\n    \"Header Char\",\n    \"Footer Char\",\n    \"Title Char\",\n    \"No Spacing\",\n    \"Subtitle Char\",\n    \"Body Text Char\",\n    \"No Spacing Char\",\n    \"Comment Char\"\n)\n\nforeach ($name in $styleNames) {\n    $style = $d.Styles($name)\n    if ($style -ne $null -and $style.Font.Name -eq $OldFont) {\n        $style.Font.Name = $NewFont\n    }\n}\n"}
